# Scheduled-runner update: refresh crafting-profit figures (currentAveragePrice,
# NQ/HQ price & profit columns H:N) across the Aegis_Profits sheets, per the
# latest market-board pull. Values only; no layout/formula/style changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3407.111
$ws.Range("I76").Value = 3194.4614
$ws.Range("J76").Value = 3960
$ws.Range("K76").Value = 3194.4614
$ws.Range("L76").Value = 3960
$ws.Range("M76").Value = -2879.4614
$ws.Range("N76").Value = -4590
$ws.Range("H79").Value = 3407.111
$ws.Range("I79").Value = 3194.4614
$ws.Range("J79").Value = 3960
$ws.Range("K79").Value = 3194.4614
$ws.Range("L79").Value = 3960
$ws.Range("M79").Value = -2102.4614
$ws.Range("N79").Value = -6144
$ws.Range("H129").Value = 921.3606600000001
$ws.Range("J129").Value = 959.9423
$ws.Range("L129").Value = 2879.8269
$ws.Range("N129").Value = -12879.8269
$ws.Range("H132").Value = 9264818
$ws.Range("I132").Value = 9621111
$ws.Range("K132").Value = 28863333
$ws.Range("M132").Value = -28860803
$ws.Range("H138").Value = 5025.22
$ws.Range("I138").Value = 1830.4736
$ws.Range("J138").Value = 6983.2905
$ws.Range("K138").Value = 5491.4208
$ws.Range("L138").Value = 20949.8715
$ws.Range("M138").Value = -351.4207999999999
$ws.Range("N138").Value = -31229.8715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31100.275
$ws.Range("I32").Value = 5762.829
$ws.Range("J32").Value = 92208.234
$ws.Range("K32").Value = 5762.829
$ws.Range("L32").Value = 92208.234
$ws.Range("M32").Value = -5475.829
$ws.Range("N32").Value = -92782.234
$ws.Range("H61").Value = 1758.3265
$ws.Range("I61").Value = 1019.37933
$ws.Range("J61").Value = 2829.8
$ws.Range("K61").Value = 1019.37933
$ws.Range("L61").Value = 2829.8
$ws.Range("M61").Value = -807.37933
$ws.Range("N61").Value = -3253.8
$ws.Range("H74").Value = 1697.4286
$ws.Range("I74").Value = 1022.875
$ws.Range("J74").Value = 1967.25
$ws.Range("K74").Value = 1022.875
$ws.Range("L74").Value = 1967.25
$ws.Range("M74").Value = -148.875
$ws.Range("N74").Value = -3715.25
$ws.Range("H77").Value = 1697.4286
$ws.Range("I77").Value = 1022.875
$ws.Range("J77").Value = 1967.25
$ws.Range("K77").Value = 5114.375
$ws.Range("L77").Value = 9836.25
$ws.Range("M77").Value = -746.375
$ws.Range("N77").Value = -18572.25
$ws.Range("H80").Value = 22795.2
$ws.Range("J80").Value = 27494
$ws.Range("L80").Value = 27494
$ws.Range("N80").Value = -29490
$ws.Range("H83").Value = 22795.2
$ws.Range("J83").Value = 27494
$ws.Range("L83").Value = 82482
$ws.Range("N83").Value = -92466
$ws.Range("H88").Value = 6500.857
$ws.Range("I88").Value = 1301.5
$ws.Range("J88").Value = 13433.333
$ws.Range("K88").Value = 1301.5
$ws.Range("L88").Value = 13433.333
$ws.Range("M88").Value = -895.5
$ws.Range("N88").Value = -14245.333
$ws.Range("H91").Value = 6500.857
$ws.Range("I91").Value = 1301.5
$ws.Range("J91").Value = 13433.333
$ws.Range("K91").Value = 1301.5
$ws.Range("L91").Value = 13433.333
$ws.Range("M91").Value = 102.5
$ws.Range("N91").Value = -16241.333
$ws.Range("H122").Value = 2199.12
$ws.Range("I122").Value = 2075.8823
$ws.Range("J122").Value = 2461
$ws.Range("K122").Value = 6227.646900000001
$ws.Range("L122").Value = 7383
$ws.Range("M122").Value = -3777.646900000001
$ws.Range("N122").Value = -12283
$ws.Range("H132").Value = 2487.1667
$ws.Range("I132").Value = 2321.5417
$ws.Range("J132").Value = 3149.6667
$ws.Range("K132").Value = 6964.625100000001
$ws.Range("L132").Value = 9449.000100000001
$ws.Range("M132").Value = -4434.625100000001
$ws.Range("N132").Value = -14509.0001
$ws.Range("H136").Value = 1758.3265
$ws.Range("I136").Value = 1019.37933
$ws.Range("J136").Value = 2829.8
$ws.Range("K136").Value = 3058.13799
$ws.Range("L136").Value = 8489.400000000001
$ws.Range("M136").Value = -508.1379900000002
$ws.Range("N136").Value = -13589.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 358.2143
$ws.Range("I22").Value = 346.9091
$ws.Range("J22").Value = 399.66666
$ws.Range("K22").Value = 346.9091
$ws.Range("L22").Value = 399.66666
$ws.Range("M22").Value = -173.9091
$ws.Range("N22").Value = -745.66666
$ws.Range("H94").Value = 36871.395
$ws.Range("I94").Value = 500599.5
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 500599.5
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -500148.5
$ws.Range("N94").Value = -2102
$ws.Range("H134").Value = 1891.3
$ws.Range("I134").Value = 2028.3043
$ws.Range("J134").Value = 1441.1428
$ws.Range("K134").Value = 6084.9129
$ws.Range("L134").Value = 4323.428400000001
$ws.Range("M134").Value = -3549.9129
$ws.Range("N134").Value = -9393.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 487.6
$ws.Range("I22").Value = 314.8
$ws.Range("J22").Value = 660.4
$ws.Range("K22").Value = 314.8
$ws.Range("L22").Value = 660.4
$ws.Range("M22").Value = 35.19999999999999
$ws.Range("N22").Value = -1360.4
$ws.Range("H31").Value = 14522.895
$ws.Range("I31").Value = 24706.309
$ws.Range("J31").Value = 1943.3823
$ws.Range("K31").Value = 24706.309
$ws.Range("L31").Value = 1943.3823
$ws.Range("M31").Value = -24411.309
$ws.Range("N31").Value = -2533.3823
$ws.Range("H34").Value = 14522.895
$ws.Range("I34").Value = 24706.309
$ws.Range("J34").Value = 1943.3823
$ws.Range("K34").Value = 24706.309
$ws.Range("L34").Value = 1943.3823
$ws.Range("M34").Value = -24504.309
$ws.Range("N34").Value = -2347.3823
$ws.Range("H132").Value = 2959.7778
$ws.Range("I132").Value = 2808.3076
$ws.Range("J132").Value = 3353.6
$ws.Range("K132").Value = 8424.9228
$ws.Range("L132").Value = 10060.8
$ws.Range("M132").Value = -5894.9228
$ws.Range("N132").Value = -15120.8
$ws.Range("H134").Value = 1911.7
$ws.Range("I134").Value = 1388
$ws.Range("K134").Value = 4164
$ws.Range("M134").Value = -1629

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 800
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 2400
$ws.Range("M22").Value = -2231
$ws.Range("H27").Value = 800
$ws.Range("I27").Value = 800
$ws.Range("K27").Value = 2400
$ws.Range("M27").Value = -2298
$ws.Range("H34").Value = 1675
$ws.Range("I34").Value = 166.66667
$ws.Range("J34").Value = 2806.25
$ws.Range("K34").Value = 500.00001
$ws.Range("L34").Value = 8418.75
$ws.Range("M34").Value = -416.00001
$ws.Range("N34").Value = -8586.75
$ws.Range("H122").Value = 50249.5
$ws.Range("I122").Value = 500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
$ws.Range("H131").Value = 1614.9302
$ws.Range("I131").Value = 547
$ws.Range("J131").Value = 1724.4615
$ws.Range("K131").Value = 1641
$ws.Range("L131").Value = 5173.3845
$ws.Range("M131").Value = 3399
$ws.Range("N131").Value = -15253.3845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 256.33334
$ws.Range("I2").Value = 285
$ws.Range("J2").Value = 231.25
$ws.Range("K2").Value = 285
$ws.Range("L2").Value = 231.25
$ws.Range("M2").Value = -172
$ws.Range("N2").Value = -457.25
$ws.Range("H102").Value = 2338
$ws.Range("I102").Value = 3536
$ws.Range("J102").Value = 1140
$ws.Range("K102").Value = 3536
$ws.Range("L102").Value = 1140
$ws.Range("M102").Value = -1914
$ws.Range("N102").Value = -4384

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4166.3335
$ws.Range("J46").Value = 4800
$ws.Range("L46").Value = 4800
$ws.Range("N46").Value = -5176

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3081.6
$ws.Range("I132").Value = 3100.9546
$ws.Range("J132").Value = 3028.375
$ws.Range("K132").Value = 9302.863799999999
$ws.Range("L132").Value = 9085.125
$ws.Range("M132").Value = -6772.863799999999
$ws.Range("N132").Value = -14145.125
